$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "height-auto" value to column L for the sub-item rows of the
# "1.11" (rows 14-21) and "1.12" (row 23) function-group sections.
# Row 22 (the "1.12" / TernaryDropdown header row) is intentionally skipped.
$rows = @(14, 15, 16, 17, 18, 19, 20, 21, 23)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 12).Value = "height-auto"
}

# Move the active selection to K24, matching the saved workbook state.
$ws.Range("K24").Select()
